$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark Texture as "Yes" for Beef Sandwich (row 3) and Chicken Sandwich (row 7),
# and mark Visuals Tested as "Yes" for Chicken Sandwich (row 7).
$ws.Range("H3").Value = "Yes"
$ws.Range("H7").Value = "Yes"
$ws.Range("J7").Value = "Yes"

# Update the active selection to match the resulting workbook state.
$ws.Range("H6").Select()
